# Realestate Update resale numbers 2025-02-26 22:41
# Append a new data row (row 91) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

# Column A: Date - force text so "2025-02-26" isn't auto-converted to a date serial.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-02-26"
$cellA.ClearFormats()

# Column B: Time (stored as plain text like the existing rows).
$ws.Cells.Item($row, 2).Value = "22:41:38"

# Column C: Weekday (plain text).
$ws.Cells.Item($row, 3).Value = "Wednesday"

# Column D: Week - force text so "08" keeps its leading zero.
$cellD = $ws.Cells.Item($row, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "08"
$cellD.ClearFormats()

# Columns E-T: numeric city values.
$ws.Cells.Item($row, 5).Value = 131207
$ws.Cells.Item($row, 6).Value = 141985
$ws.Cells.Item($row, 7).Value = 173118
$ws.Cells.Item($row, 8).Value = 160254
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 146780
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 193954
$ws.Cells.Item($row, 14).Value = 115467
$ws.Cells.Item($row, 15).Value = 46809
$ws.Cells.Item($row, 16).Value = 29556
$ws.Cells.Item($row, 17).Value = 69669
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 50187
$ws.Cells.Item($row, 20).Value = -1
